$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-27 07:57:42"
$wsZhCn.Range("G3").Value = "2016-01-27 07:58:25"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-27 07:57:55"
$wsDeDe.Range("G3").Value = "2016-01-27 07:58:47"
